$d = $word.ActiveDocument

# The edit removes three whole paragraphs that sit between
# "5.4 Estrutura exigida" (kept, with its own image) and
# "5.7 Manutenção" (kept):
#   - "5.5 Abrangência" (heading text + line break + its image)
#   - the blank paragraph right after it
#   - "5.6 integração"
#
# Locate the boundary paragraphs by their text instead of a fixed index so
# the script keeps working even if paragraph numbering shifts.

$startIndex = 0
$endIndex = 0
$count = $d.Paragraphs.Count

for ($i = 1; $i -le $count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($startIndex -eq 0 -and $text -like "5.5 Abrang*") {
        $startIndex = $i
    }
    if ($startIndex -ne 0 -and $text -like "5.6*") {
        $endIndex = $i
        break
    }
}

if ($startIndex -ne 0 -and $endIndex -ne 0) {
    $startPara = $d.Paragraphs.Item($startIndex)
    $endPara = $d.Paragraphs.Item($endIndex)
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
